# Daily attendance update - 2025-08-28
# Fill in the AF (2025-08-27, "India Holiday") and AG (2025-08-28, "WFO"/"WFH")
# columns for every employee row on the August 2025 tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WCS_Team_August_2025")
$ws.Activate()

# Row -> attendance status for 2025-08-28 (column AG)
$attendance = @{
    3  = "WFH"
    4  = "WFO"
    5  = "WFO"
    6  = "WFH"
    7  = "WFH"
    8  = "WFH"
    9  = "WFH"
    10 = "WFH"
    11 = "WFO"
    12 = "WFH"
    13 = "WFH"
    14 = "WFH"
    15 = "WFH"
    16 = "WFH"
    17 = "WFH"
    18 = "WFH"
}

foreach ($row in 3..18) {
    # Column AF (2025-08-27) was a company holiday for everyone.
    $ws.Range("AF$row").Value = "India Holiday"

    # Column AG (2025-08-28): copy the conditional-format style from an
    # existing cell that already carries the right WFO/WFH look, then set
    # the value - this reuses the workbook's existing "Good"/"Neutral"
    # cell styles instead of creating new ones. I3 is a stable "WFO"-styled
    # cell and F3 is a stable "WFH"-styled cell; neither is touched by this
    # script, so they remain valid copy sources for every loop iteration.
    $status = $attendance[$row]
    if ($status -eq "WFO") {
        $ws.Range("I3").Copy()
    } else {
        $ws.Range("F3").Copy()
    }
    $ws.Range("AG$row").PasteSpecial(-4122)
    $ws.Range("AG$row").Value = $status
}

$excel.CutCopyMode = 0

# Leave the selection where the last edit happened.
$ws.Range("AG11").Select()
